$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (updated electricity spot prices for the new day)
$ws.Range("A2").Value = 46060
$ws.Range("B2").Value = 0.29
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.05
$ws.Range("E2").Value = -0.14
$ws.Range("F2").Value = -0.16
$ws.Range("G2").Value = -0.17
$ws.Range("H2").Value = -0.11
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -0.04
$ws.Range("M2").Value = -0.1
$ws.Range("N2").Value = -0.05
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.01
$ws.Range("Q2").Value = 0.17
$ws.Range("R2").Value = 0.59
$ws.Range("S2").Value = 0.17
$ws.Range("T2").Value = 0.44
$ws.Range("U2").Value = 1.55
$ws.Range("V2").Value = 3.09
$ws.Range("W2").Value = 3.81
$ws.Range("X2").Value = 2.24
$ws.Range("Y2").Value = 1.44
$ws.Range("Z2").Value = 0.54
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 2.64
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 3.45
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 1.84
$ws.Range("AG2").Value = "0h-18h"
